$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, pushing existing rows 68-74 down to 69-75.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly record.
$ws.Range("A68").Value2 = 2
$ws.Range("B68").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C68").Value2 = "Coquimbo"
$ws.Range("D68").Value2 = 44615
$ws.Range("E68").Value2 = 4
$ws.Range("F68").Value2 = 100112030
$ws.Range("G68").Value2 = "Poroto granado"
$ws.Range("H68").Value2 = "Sin especificar"
$ws.Range("I68").Value2 = "Primera"
$ws.Range("J68").Value2 = 500
$ws.Range("K68").Value2 = 21000
$ws.Range("L68").Value2 = 23000
$ws.Range("M68").Value2 = 22000
$ws.Range("N68").Value2 = "$/malla 25 kilos"
$ws.Range("O68").Value2 = "Provincia de Limarí"
$ws.Range("P68").Value2 = 880
$ws.Range("Q68").Value2 = 25
$ws.Range("R68").Value2 = "Hortaliza"
